$d = $word.ActiveDocument

# The target bullet is the last paragraph in the document:
#   "Create installer <Wingdings arrow> Gen report"
# It becomes:
#   "Create installer, checking installer <Wingdings arrow> Gen report "
$para = $d.Paragraphs.Last

# 1) Turn "Create installer " into "Create installer, checking installer "
#    (adds ", checking installer" right after "Create installer", keeping
#    the single space that was already separating it from the Wingdings
#    arrow symbol run).
$r = $para.Range
$r.Find.Execute("Create installer ", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Create installer, checking installer ", 2)

# 2) Append a trailing space after " Gen report", right before the
#    paragraph mark, as its own run.
$para2 = $d.Paragraphs.Last
$tail = $para2.Range
$tail.Collapse(0)
$tail.MoveEnd(1, -1)
$tail.InsertAfter(" ")
